# edit.ps1 - apply "slides for lesson 06" changes to lesson_05_ejb.pptx
#
# Summary of changes:
#  1. Insert a new slide at position 8 titled "@PostConstruct" (the old slide 8,
#     "Container Deployment", and everything after it, shifts down by one position).
#     We do this by duplicating the old slide 8 (so the duplicate keeps the
#     "Container Deployment" content at the new position 9) and then rewriting the
#     original slide 8 in place with the new "@PostConstruct" content.
#  2. Slide 1 (title slide): merge 3 separate runs "Lesson " / "05: " / "EJB" into a
#     single run "Lesson 05: EJB".
#  3. Slide 13 (was slide 12 before the insertion, "Git Repository Modules"): merge 4
#     separate runs "Exercises " / "for Lesson " / "05 " / "(see documentation)" into a
#     single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: duplicate slide 8 ("Container Deployment") so that content survives,
# unchanged, at the new position 9.
# ---------------------------------------------------------------------------
$srcSlide = $p.Slides.Item(8)
$srcSlide.Duplicate() | Out-Null

# ---------------------------------------------------------------------------
# Step 2: rewrite slide 8 in place as the new "@PostConstruct" slide.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

# --- Title -------------------------------------------------------------
$titleShape = $slide8.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "@PostConstruct"
$atSign = $titleRange.Characters(1, 1)
$atSign.Text = $atSign.Text
$rest = $titleRange.Characters(2, 13)
$rest.Text = $rest.Text

# --- Content placeholder: resize/reposition -----------------------------
$bodyShape = $slide8.Shapes.Item(2)
$bodyShape.Left = (288324 / 12700.0)
$bodyShape.Top = (1825624 / 12700.0)
$bodyShape.Width = (11788346 / 12700.0)
$bodyShape.Height = (4863499 / 12700.0)
$bodyShape.TextFrame.AutoSize = 2

# --- Content placeholder: body text --------------------------------------
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "x"

$q1 = [char]8220
$q2 = [char]8221

$para1 = "The Container, before doing dependency injection, needs to create an instance of the EJB with " + $q1 + "new" + $q2
$para2 = "This means that the code of the constructor is called BEFORE dependency injection (DI) is done"
$para3 = "If you need to access an injected variable in the constructor, you will hence get a null pointer exception"
$para4 = "A method marked with @PostConstruct will be executed AFTER the constructor and DI "
$para5 = "so, useful when you need initializing code relying on injected variables"

$fullBody = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5
$bodyRange.Text = $fullBody

# paragraph 2: italicize the word "constructor"
$p2Start = $fullBody.IndexOf($para2) + 1
$ctorIdx = $para2.IndexOf("constructor")
$ctorRange = $bodyRange.Characters($p2Start + $ctorIdx, 11)
$ctorRange.Font.Italic = -1

# paragraph 4: split into "A method marked with " / "@" / "PostConstruct" / " will be executed..."
$p4Start = $fullBody.IndexOf($para4) + 1
$atIdx = $para4.IndexOf("@")
$atRange = $bodyRange.Characters($p4Start + $atIdx, 1)
$atRange.Text = $atRange.Text
$postConstructRange = $bodyRange.Characters($p4Start + $atIdx + 1, 13)
$postConstructRange.Text = $postConstructRange.Text

# paragraph 5: indent level 2 (lvl="1"), split into "s" / "o, useful..." / "injected variables"
$lastPara = $bodyRange.Paragraphs(5, 1)
$lastPara.IndentLevel = 2
$p5Start = $fullBody.IndexOf($para5) + 1
$sRange = $bodyRange.Characters($p5Start, 1)
$sRange.Text = $sRange.Text
$ivIdx = $para5.IndexOf("injected variables")
$ivRange = $bodyRange.Characters($p5Start + $ivIdx, 19)
$ivRange.Text = $ivRange.Text

# ---------------------------------------------------------------------------
# Step 3: slide 1 title - merge "Lesson " + "05: " + "EJB" into one run.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleRange1 = $slide1.Shapes.Item(1).TextFrame.TextRange
$t1 = $titleRange1.Text
$lessonIdx = $t1.IndexOf("Lesson")
$lessonLen = $t1.Length - $lessonIdx
$lessonRange = $titleRange1.Characters($lessonIdx + 1, $lessonLen)
$lessonRange.Text = "Lesson 05: EJB"

# ---------------------------------------------------------------------------
# Step 4: slide 13 (Git Repository Modules, shifted from 12 -> 13) - merge
# "Exercises " + "for Lesson " + "05 " + "(see documentation)" into one run.
# ---------------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$bodyRange13 = $slide13.Shapes.Item(2).TextFrame.TextRange
$t13 = $bodyRange13.Text
$exIdx = $t13.IndexOf("Exercises")
$exLen = $t13.Length - $exIdx
$exRange = $bodyRange13.Characters($exIdx + 1, $exLen)
$exRange.Text = "Exercises for Lesson 05 (see documentation)"
